# Add a new row of config data to Sheet1 (row 16), as requested in client feedback.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D16").Value = "shareholdings_summary"
$ws.Range("A16").Value = "total_equity_shares"
$ws.Range("C16").Value = "number_issued_under_share_structure"
$ws.Range("B16").Value = "Single"
$ws.Range("E16").Value = "total_equity_shares"

# Match the formatting of the row above (A15 uses the wrap-text, no-border style)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("A16").Select()
